$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new transaction data
$ws.Range("A2").Value = "Переводы"
$ws.Range("B2").Value = 54
$ws.Range("C2").Value = "29/9/2023"

# Remove row 3 entirely (the second "Заработная плата" entry)
$ws.Rows("3:3").Delete()
